$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph -
# it sits between a blank paragraph and the "© 2020 ..." copyright
# paragraph. All three paragraphs (blank, "Ver no Jupiter...", and the
# copyright line) must be removed, leaving the "LOQ4057: ..." paragraph
# followed directly by the blank paragraph that used to precede the
# page-break paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Ver no Jupiter Salvar em pdf Salvar em docx") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $prevPara = $target.Previous()
    $nextPara = $target.Next()
    $start = $prevPara.Range.Start
    $end = $nextPara.Range.End
    $d.Range($start, $end).Delete()
}
